# Update the header row of the active sheet ("s2") with new column labels
# and move the selection, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "column_1"
$ws.Range("B1").Value = "column_2"

[void]$ws.Range("D16").Select()
